$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'71.150.04"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +6.18%  "
$ws.Range("D3").Value = "'3.663.75"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +5.88%  "
$ws.Range("D5").Value = "'596.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.99%  "
$ws.Range("D6").Value = "'195.70"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.58%  "
$ws.Range("D7").Value = "'0.650"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.62%  "
$ws.Range("D8").Value = "'3.656.46"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.86%  "
$ws.Range("E9").Value = "  -0.06%  "
$ws.Range("E10").Value = "  +5.50%  "
$ws.Range("E11").Value = "  +4.96%  "
$ws.Range("D12").Value = "'58.98"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.88%  "
$ws.Range("E13").Value = "  +6.27%  "
$ws.Range("D14").Value = "'10.01"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.59%  "
$ws.Range("D15").Value = "'4.252.05"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.77%  "
$ws.Range("D16").Value = "'19.97"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +6.79%  "
$ws.Range("D17").Value = "'3.665.08"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.87%  "
$ws.Range("D18").Value = "'71.147.56"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +6.25%  "
$ws.Range("D19").Value = "'12.87"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +6.25%  "
$ws.Range("E20").Value = "  +3.19%  "
$ws.Range("E21").Value = "  +5.48%  "
$ws.Range("D22").Value = "'494.19"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.68%  "
$ws.Range("D23").Value = "'19.08"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +15.23%  "
$ws.Range("E24").Value = "  +0.45%  "
$ws.Range("E25").Value = "  +2.23%  "
$ws.Range("D26").Value = "'91.88"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.45%  "
$ws.Range("D27").Value = "'3.19"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +8.67%  "
$ws.Range("D28").Value = "'11.59"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.54%  "
$ws.Range("D29").Value = "'9.70"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +7.28%  "
$ws.Range("D30").Value = "'33.09"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.69%  "
$ws.Range("D31").Value = "'7.93"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +11.27%  "
$ws.Range("E32").Value = "  +9.67%  "
$ws.Range("D33").Value = "'633.64"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.68%  "
$ws.Range("D34").Value = "'12.33"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.48%  "
$ws.Range("D35").Value = "'65.97"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.40%  "
$ws.Range("D36").Value = "'40.96"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +12.38%  "
$ws.Range("D37").Value = "'0.0₃0840"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +11.52%  "
$ws.Range("D38").Value = "'0.416"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +8.57%  "
$ws.Range("E39").Value = "  -0.91%  "
$ws.Range("E41").Value = "  +2.34%  "
$ws.Range("D42").Value = "'3.334.26"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.09%  "
$ws.Range("D43").Value = "'3.18"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +9.48%  "
$ws.Range("D44").Value = "'2.88"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +14.83%  "
$ws.Range("E45").Value = "  +6.52%  "
$ws.Range("D46").Value = "'2.91"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.83%  "
$ws.Range("E47").Value = "  +3.53%  "
$ws.Range("D48").Value = "'3.30"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.06%  "
$ws.Range("D49").Value = "'9.34"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +7.21%  "
$ws.Range("E50").Value = "  +1.67%  "
$ws.Range("E51").Value = "  +0.05%  "
